# Final sap flow upscaling analysis
# Corrects mislabeled "Irrigation stop" entries in column E (Treatment) to
# "Control" for a set of rows on Sheet1, and leaves the final selection on
# the last edited block (E459:E463), matching the author's last action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contiguous row ranges (1-based worksheet rows) in column E whose value
# changes from "Irrigation stop" to "Control".
$ranges = @(
    "E77:E79",
    "E81:E81",
    "E118:E153",
    "E328:E347",
    "E356:E356",
    "E387:E391",
    "E393:E404",
    "E423:E436",
    "E459:E463"
)

foreach ($rangeAddress in $ranges) {
    $ws.Range($rangeAddress).Value = "Control"
}

# Reflect the final cell selection left by the author after the last edit.
$ws.Range("E459:E463").Select()
